$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 (August) label to reflect new "through" date
$ws.Range("A9").Value = "August (through 08-27)"

# Update row 9 (August) yearly totals
$ws.Range("B9").Value = 28
$ws.Range("C9").Value = 67
$ws.Range("D9").Value = 78
$ws.Range("E9").Value = 57
$ws.Range("G9").Value = 152
$ws.Range("H9").Value = 136

# Update row 10 (Total) yearly totals
$ws.Range("B10").Value = 190
$ws.Range("C10").Value = 369
$ws.Range("D10").Value = 543
$ws.Range("E10").Value = 482
$ws.Range("G10").Value = 773
$ws.Range("H10").Value = 1050
